$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update row 2 data values (new test user: Sandeep) ---
$ws.Range("A2").Value = "Sandeep"
$ws.Range("B2").Value = "sandeep@gmail.com"
$ws.Range("C2").Value = 7817008251
$ws.Range("D2").Value = "Vadodara"
$ws.Range("F2").Value = "Sandeep@123"
$ws.Range("G2").Value = "Sandeep@123"

# --- Update hyperlink targets to match the new email / pin values ---
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.Address = "mailto:sandeep@gmail.com"
    } elseif ($addr -eq '$F$2') {
        $h.Address = "mailto:Sandeep@123"
    } elseif ($addr -eq '$G$2') {
        $h.Address = "mailto:Sandeep@123"
    }
}

# --- Selection moved to H1 ---
[void]$ws.Range("H1").Select()

# --- Window size/position changed (maximized) ---
$excel.ActiveWindow.Top = -108
$excel.ActiveWindow.Left = -108
$excel.ActiveWindow.Width = 23256
$excel.ActiveWindow.Height = 12456
